$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets, LEAGUE_RULES and GAMES, at the end (after
#    TEAM_PLAYER_MAPPINGS), preserving tab order / TEAM_PLAYER_MAPPINGS as
#    the active sheet.
# ---------------------------------------------------------------------------
$wsTeamNames   = $wb.Worksheets.Item("TEAMNAMES")
$wsTeamPlayers = $wb.Worksheets.Item("TEAM_PLAYER_MAPPINGS")

$wsRules = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTeamPlayers)
$wsRules.Name = "LEAGUE_RULES"

$wsGames = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsRules)
$wsGames.Name = "GAMES"

# ---------------------------------------------------------------------------
# 2. TEAM_PLAYER_MAPPINGS: insert a new UNIQUE_ID column between LAST_NAME
#    and PLAYER_ROLE. Done by re-writing the cell contents directly (rather
#    than a real column insert) so column-width metadata on column C stays
#    put on column C.
# ---------------------------------------------------------------------------
$playerRows = @(
    @("FIRST_NAME",   "LAST_NAME",   "UNIQUE_ID", "PLAYER_ROLE",   "TEAM_INITIALS"),
    @("Rohit",        "Sharma",      "RS",        "Batsman",       "MI"),
    @("Harbhajan",    "Singh",       "HS",        "Bowler",        "MI"),
    @("Lasith",       "Malinga",     "LM",        "Bowler",        "MI"),
    @("Kieron",       "Pollard",     "KP",        "All-Rounder",   "MI"),
    @("MS",           "Dhoni",       "MSD",       "WicketKeeper",  "CSK"),
    @("Ravindra",     "Jadeja",      "RJ",        "All-Rounder",   "CSK"),
    @("Dwayne ",      "Bravo",       "DB",        "All-Rounder",   "CSK"),
    @("Virat",        "Kohli",       "VK",        "Batsman",       "RCB"),
    @("AB",           "deVilliers",  "ABdeV",     "WicketKeeper",  "RCB"),
    @("Shane",        "Watson",      "SW",        "All-Rounder",   "RCB"),
    @("KL",           "Rahul",       "KLR",       "Batsman",       "RCB"),
    @("Gautam",       "Gambhir",     "GG",        "Batsman",       "KKR"),
    @("Yusuf",        "Pathan",      "YP",        "All-Rounder",   "KKR"),
    @("Robin",        "Uthappa",     "RU",        "WicketKeeper",  "KKR"),
    @("Sunil",        "Narine",      "SN",        "Bowler",        "KKR"),
    @("David",        "Warner",      "DW",        "Batsman",       "SRH"),
    @("Shikhar",      "Dhawan",      "SD",        "Batsman",       "SRH"),
    @("Bhvaneshwar",  "Kumar",       "BK",        "Bowler",        "SRH"),
    @("Ashish",       "Nehra",       "AN",        "Bowler",        "SRH")
)

for ($i = 0; $i -lt $playerRows.Length; $i++) {
    $r = $i + 1
    $row = $playerRows[$i]
    $wsTeamPlayers.Cells.Item($r, 3).Value = $row[2]
    $wsTeamPlayers.Cells.Item($r, 4).Value = $row[3]
    $wsTeamPlayers.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# 3. LEAGUE_RULES sheet content
# ---------------------------------------------------------------------------
$wsRules.Range("A1").Value = "METRIC"
$wsRules.Range("B1").Value = "OPERATOR"
$wsRules.Range("C1").Value = "SCORE"
$wsRules.Range("D1").Value = "POINTS"
$wsRules.Range("E1").Value = "RATIO"

$wsRules.Range("A2").Value = "RUNS"
$wsRules.Range("B2").Value = "EQUALS"
$wsRules.Range("C2").Value = 0
$wsRules.Range("D2").Value = -5
$wsRules.Range("E2").Value = "TOTAL"

$wsRules.Range("A3").Value = "RUNS"
$wsRules.Range("B3").Value = "GREATER_THAN"
$wsRules.Range("C3").Value = 50
$wsRules.Range("D3").Value = 5
$wsRules.Range("E3").Value = "PER_SCORE"

# Row 4: only C4 has content -- a stray quote-prefixed, value-less cell.
$wsRules.Range("C4").Value = "'"
$wsRules.Range("C4").Value = ""

$wsRules.Range("A16").Value = "RUNS"
$wsRules.Range("B16").Value = "LESS_THAN"
$wsRules.Range("E16").Value = "PER_SCORE"

$wsRules.Range("A17").Value = "WICKETS"
$wsRules.Range("B17").Value = "GREATER_THAN"
$wsRules.Range("E17").Value = "TOTAL"

$wsRules.Range("A18").Value = "CATCHES"
$wsRules.Range("B18").Value = "EQUALS"

$wsRules.Range("A1,A2:A14").Validation.Add(3, 1, 1, "=`$A`$16:`$A`$18")
$wsRules.Range("B2:B15").Validation.Add(3, 1, 1, "=`$B`$16:`$B`$18")
$wsRules.Range("E2:E15").Validation.Add(3, 1, 1, "=`$E`$16:`$E`$17")

$wsRules.Columns("B:B").ColumnWidth = 13.333333333333334
$wsRules.Columns("C:C").ColumnWidth = 12.666666666666666
$wsRules.Columns("D:D").ColumnWidth = 13.333333333333334
$wsRules.Columns("E:E").ColumnWidth = 15.833333333333334

# ---------------------------------------------------------------------------
# 4. GAMES sheet content
# ---------------------------------------------------------------------------
$wsGames.Range("A1").Value = "TEAM1"
$wsGames.Range("B1").Value = "TEAM2"
$wsGames.Range("C1").Value = "DATE( MM/DD/YYYY)"
$wsGames.Range("D1").Value = "TIME(hh:mm:ss)"
$wsGames.Range("E1").Value = "VENUE"

$wsGames.Range("A2").Value = "KKR"
$wsGames.Range("B2").Value = "MI"
$wsGames.Range("C2").Value = (Get-Date -Year 2017 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0)
$wsGames.Range("C2").NumberFormat = "mm-dd-yy"
$wsGames.Range("D2").Value = 0.66666666666666663
$wsGames.Range("D2").NumberFormat = "h:mm:ss"
$wsGames.Range("E2").Value = "MUMBAI"

$wsGames.Range("A3").Value = "CSK"
$wsGames.Range("B3").Value = "RCB"
$wsGames.Range("C3").Value = (Get-Date -Year 2017 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0)
$wsGames.Range("C3").NumberFormat = "mm-dd-yy"
$wsGames.Range("D3").Value = 0.83333333333333337
$wsGames.Range("D3").NumberFormat = "h:mm:ss"
$wsGames.Range("E3").Value = "MUMBAI"

$wsGames.Columns("C:C").ColumnWidth = 19
$wsGames.Columns("D:D").ColumnWidth = 17.333333333333332

# ---------------------------------------------------------------------------
# 5. Selections: restore per-sheet selection state, keep TEAM_PLAYER_MAPPINGS
#    as the active/tab-selected sheet.
# ---------------------------------------------------------------------------
$wsTeamNames.Range("A29").Select()
$wsRules.Range("B4").Select()
$wsGames.Range("D3").Select()

$wsTeamPlayers.Select()
$wsTeamPlayers.Range("G11").Select()
